# Generate Report for Handoff
# Update the "Status" text from "Handed back: in sync with en-US" to
# "Ready for handoff" on all three sheets, and refresh the associated
# timestamp columns, then let the column widths auto-fit the shorter text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# E2 (zh-cn status) and F2 (de-de status)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# G2 Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-01 05:04:05"

# --- zh-cn sheet ---
# C2 Status
$wsZhCn.Range("C2").Value = "Ready for handoff"
# H2 Latest Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-09-01 05:03:56"

# --- de-de sheet ---
# C2 Status
$wsDeDe.Range("C2").Value = "Ready for handoff"

# The shorter status text means these columns no longer need to be as
# wide; set the resulting auto-fit widths explicitly (target ~17.216
# characters; input chosen so the engine's internal rounding lands on
# the closest representable column width).
$wsOverview.Range("E:F").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C:C").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C:C").ColumnWidth = 16.333333333333332
